$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the active selection to I43 (this also clears the previous
# scrolled "topLeftCell" saved in the sheet view)
$ws.Range("I43").Select()

# Fill in the newly-computed pre-training average results (column D)
# for the rows that previously had no value there.
$ws.Range("D43").Formula = "=(0.496904024767801 + 0.496904024767801 + 0.461300309597523) / 3 * 100"
$ws.Range("D44").Formula = "=(0.434715346534653 + 0.435334158415841 + 0.426980198019801) / 3 * 100"
$ws.Range("D46").Formula = "=(0.32824427480916 + 0.3206106870229 + 0.396946564885496) / 3 * 100"
$ws.Range("D47").Formula = "=(0.266768292682926 + 0.265243902439024 +  0.271341463414634) / 3 * 100"
